$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.167.79'
$ws.Range('E2').Value = '  -2.90%  '
$ws.Range('D3').Value = '1.848.84'
$ws.Range('E3').Value = '  -1.93%  '
$ws.Range('E4').Value = '  -0.18%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7041'
$ws.Range('E5').Value = '  -4.31%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '238.58'
$ws.Range('E6').Value = '  -1.64%  '
$ws.Range('E7').Value = '  -0.30%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3052'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07457'
$ws.Range('E9').Value = '  +4.00%  '
$ws.Range('E10').Value = '  -4.85%  '
$ws.Range('E11').Value = '  -2.54%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.7266'
$ws.Range('E12').Value = '  -3.75%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.843.72'
$ws.Range('E13').Value = '  -7.44%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.237'
$ws.Range('E14').Value = '  -2.92%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '88.78'
$ws.Range('E15').Value = '  -4.10%  '
$ws.Range('D16').Value = '29.077.00'
$ws.Range('E16').Value = '  -3.33%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '5.764'
$ws.Range('E17').Value = '  -5.98%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '238.72'
$ws.Range('E18').Value = '  -4.27%  '
$ws.Range('E19').Value = '  -3.38%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007618'
$ws.Range('E20').Value = '  -3.01%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.9988'
$ws.Range('E21').Value = '  -0.32%  '
$ws.Range('D22').Value = '2.081.95'
$ws.Range('E22').Value = '  -3.76%  '
$ws.Range('E23').Value = '  -0.16%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '7.601'
$ws.Range('E24').Value = '  -4.01%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '8.996'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '161.06'
$ws.Range('E26').Value = '  -1.92%  '
$ws.Range('E27').Value = '  -7.15%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.08'
$ws.Range('E28').Value = '  -3.06%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.973'
$ws.Range('E29').Value = '  -3.52%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.397'
$ws.Range('E30').Value = '  -5.57%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.548'
$ws.Range('E31').Value = '  -0.15%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.490'
$ws.Range('E32').Value = '  -2.99%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.982'
$ws.Range('E33').Value = '  -4.66%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05165'
$ws.Range('E34').Value = '  -3.04%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.186'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.031'
$ws.Range('E36').Value = '  +3.53%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.7035'
$ws.Range('E37').Value = '  -8.44%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.661'
$ws.Range('E38').Value = '  -2.49%  '
$ws.Range('E39').Value = '  -4.42%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.682'
$ws.Range('E40').Value = '  -2.78%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9434'
$ws.Range('E41').Value = '  +7.55%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.010'
$ws.Range('E42').Value = '  -0.56%  '
$ws.Range('D43').Value = '1.074.44'
$ws.Range('E43').Value = '  -2.43%  '
$ws.Range('E44').Value = '  -5.82%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '70.04'
$ws.Range('E45').Value = '  -3.29%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.000'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '102.71'
$ws.Range('E47').Value = '  -1.54%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.739'
$ws.Range('E48').Value = '  -6.09%  '
$ws.Range('D49').Value = '1.988.47'
$ws.Range('E49').Value = '  -4.35%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.052'
$ws.Range('E50').Value = '  -6.63%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '9.134'
$ws.Range('E51').Value = '  -4.45%  '
